# Generate Report for Handoff
# Update the GUID-based file name and timestamps across the Overview, zh-cn
# and de-de sheets to reflect a newly generated handoff report.

$wb = $excel.ActiveWorkbook

$oldGuid = "a993b9c8-f774-4f2d-b45c-3116fcfd1fe8"
$newGuid = "ae52835f-d048-4cc5-9112-25f87fdbc015"

$oldHash = "e67de063da5112b2ce7b0ccce78a0c935b605417"
$newHash = "99febb098e46c834adcf79e00bd86b280ae57a6d"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-26-18 07:26:40"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-18 07:26:37"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-18 07:26:40"
